$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "58.069.02"
Set-TextValue $ws.Range("E2") "  -3.90%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.128.17"
Set-TextValue $ws.Range("E3") "  -5.49%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.05%  "

# Row 5
Set-TextValue $ws.Range("D5") "520.72"
Set-TextValue $ws.Range("E5") "  -6.64%  "

# Row 6
Set-TextValue $ws.Range("D6") "134.42"
Set-TextValue $ws.Range("E6") "  -5.33%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -0.10%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.130.04"
Set-TextValue $ws.Range("E8") "  -5.46%  "

# Row 9
Set-TextValue $ws.Range("E9") "  -5.77%  "

# Row 10
Set-TextValue $ws.Range("D10") "7.24"
Set-TextValue $ws.Range("E10") "  -7.86%  "

# Row 11
Set-TextValue $ws.Range("E11") "  -9.02%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.380"
Set-TextValue $ws.Range("E12") "  -6.75%  "

# Row 13
Set-TextValue $ws.Range("D13") "3.660.31"
Set-TextValue $ws.Range("E13") "  -5.61%  "

# Row 14
Set-TextValue $ws.Range("E14") "  -1.63%  "

# Row 15
Set-TextValue $ws.Range("D15") "25.44"

# Row 16
Set-TextValue $ws.Range("D16") "3.125.51"
Set-TextValue $ws.Range("E16") "  -5.84%  "

# Row 17
Set-TextValue $ws.Range("D17") "57.994.57"
Set-TextValue $ws.Range("E17") "  -4.04%  "

# Row 18
Set-TextValue $ws.Range("E18") "  -8.33%  "

# Row 19
Set-TextValue $ws.Range("E19") "  -5.38%  "

# Row 20
Set-TextValue $ws.Range("D20") "12.96"
Set-TextValue $ws.Range("E20") "  -9.75%  "

# Row 21
Set-TextValue $ws.Range("D21") "7.94"
Set-TextValue $ws.Range("E21") "  -8.42%  "

# Row 22
Set-TextValue $ws.Range("D22") "343.15"
Set-TextValue $ws.Range("E22") "  -8.26%  "

# Row 23
Set-TextValue $ws.Range("E23") "  +0.02%  "

# Row 24
Set-TextValue $ws.Range("D24") "68.80"

# Row 25
Set-TextValue $ws.Range("D25") "0.506"
Set-TextValue $ws.Range("E25") "  -6.13%  "

# Row 26
Set-TextValue $ws.Range("D26") "3.250.49"
Set-TextValue $ws.Range("E26") "  -5.65%  "

# Row 27
Set-TextValue $ws.Range("B27") "Kaspa"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D27") "0.166"
Set-TextValue $ws.Range("E27") "  -3.50%  "

# Row 28
Set-TextValue $ws.Range("B28") "PEPE"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D28") "0.0₃0948"
Set-TextValue $ws.Range("E28") "  -7.27%  "

# Row 29
Set-TextValue $ws.Range("D29") "0.996"
Set-TextValue $ws.Range("E29") "  +0.15%  "

# Row 30
Set-TextValue $ws.Range("E30") "  +0.02%  "

# Row 31
Set-TextValue $ws.Range("D31") "6.71"
Set-TextValue $ws.Range("E31") "  -6.48%  "

# Row 32
Set-TextValue $ws.Range("E32") "  -9.58%  "

# Row 33
Set-TextValue $ws.Range("D33") "21.53"
Set-TextValue $ws.Range("E33") "  -4.75%  "

# Row 34
Set-TextValue $ws.Range("D34") "6.81"
Set-TextValue $ws.Range("E34") "  -10.41%  "

# Row 35
Set-TextValue $ws.Range("E35") "  -2.16%  "

# Row 36
Set-TextValue $ws.Range("D36") "157.77"
Set-TextValue $ws.Range("E36") "  -5.02%  "

# Row 37
Set-TextValue $ws.Range("D37") "4.75"
Set-TextValue $ws.Range("E37") "  -7.93%  "

# Row 38
Set-TextValue $ws.Range("D38") "6.16"
Set-TextValue $ws.Range("E38") "  -8.17%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.35"
Set-TextValue $ws.Range("E39") "  -11.34%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0689"
Set-TextValue $ws.Range("E40") "  -6.02%  "

# Row 41
Set-TextValue $ws.Range("D41") "3.158.92"
Set-TextValue $ws.Range("E41") "  -5.43%  "

# Row 42
Set-TextValue $ws.Range("D42") "40.45"
Set-TextValue $ws.Range("E42") "  -3.79%  "

# Row 43
Set-TextValue $ws.Range("D43") "24.06"
Set-TextValue $ws.Range("E43") "  -10.60%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.695"
Set-TextValue $ws.Range("E44") "  -7.76%  "

# Row 45
Set-TextValue $ws.Range("E45") "  -3.58%  "

# Row 46
Set-TextValue $ws.Range("D46") "3.89"
Set-TextValue $ws.Range("E46") "  -6.76%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.00"
Set-TextValue $ws.Range("E47") "  -0.03%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.43"
Set-TextValue $ws.Range("E48") "  -9.55%  "

# Row 49
Set-TextValue $ws.Range("D49") "2.255.93"
Set-TextValue $ws.Range("E49") "  -4.99%  "

# Row 50
Set-TextValue $ws.Range("D50") "6.17"
Set-TextValue $ws.Range("E50") "  -4.39%  "

# Row 51
Set-TextValue $ws.Range("E51") "  -4.54%  "
